# Daily attendance processing - reorder "Recorded By" entries in column G
# so that an exact "System" token is moved to the end of the comma-separated
# list (case-sensitive). If there is no exact "System" token but there are
# multiple entries, the order of the entries is reversed.

# NOTE: the -ceq / -cne / -clike / -cmatch "case sensitive" operators are not
# actually case-sensitive in this runtime, so a manual char-code comparison
# is used instead to distinguish "System" from "system".
function Test-CaseSensitiveEquals($s1, $s2) {
    if ($s1.Length -ne $s2.Length) { return $false }
    $c1 = $s1.ToCharArray()
    $c2 = $s2.ToCharArray()
    for ($i = 0; $i -lt $c1.Count; $i++) {
        if ([int]$c1[$i] -ne [int]$c2[$i]) { return $false }
    }
    return $true
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if (-not ($val -is [string])) { continue }
    if ($val.IndexOf(",") -lt 0) { continue }

    $parts = $val -split ", "

    if ($parts.Count -le 1) { continue }

    $hasSystem = $false
    foreach ($p in $parts) {
        if (Test-CaseSensitiveEquals $p "System") { $hasSystem = $true }
    }

    if ($hasSystem) {
        $newParts = @()
        $removedOne = $false
        foreach ($p in $parts) {
            if ((-not $removedOne) -and (Test-CaseSensitiveEquals $p "System")) {
                $removedOne = $true
            } else {
                $newParts += $p
            }
        }
        $newParts += "System"
    } else {
        $newParts = @()
        for ($i = $parts.Count - 1; $i -ge 0; $i--) { $newParts += $parts[$i] }
    }

    $newVal = [string]::Join(", ", $newParts)

    if ($newVal -ne $val) {
        $cell.Value2 = $newVal
    }
}
